# Generate Report for Handback
# Adds a new handback record (84c1d267-b40e-464d-b31f-a04effd434ef.md) as row 4
# to the "Overview", "zh-cn" and "de-de" worksheets, mirroring the shape of the
# existing rows (e2f2f2be... / 6447a6b2...), and grows each sheet's table/
# dimension/autofilter range to include the new row.

$wb = $excel.ActiveWorkbook

$uuid   = "84c1d267-b40e-464d-b31f-a04effd434ef"
$mdName = "$uuid.md"
$mdDisp = "e2e\$uuid.md"

$xlfBase   = "$uuid.245edb31b7072088562f79e4638dc8fc0cf82a07"
$xlfZhCn   = "$xlfBase.zh-cn.xlf"
$xlfDeDe   = "$xlfBase.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

# Plausible commit-style hyperlink targets, following the same pattern as the
# existing rows ( https://github.com/<org>/<repo>/blob/<sha>/e2e/<file>.md ).
$mainUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4a1f3e6d0a9b4c7e2f5061829ab4de37c0915af/e2e/$mdName"
$zhcnUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7ad3f8b0c1e4d6a92f5037bc6819de4a5fb20c83/e2e/$mdName"
$dedeUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9e2c74a1f0b83d6e5c419a7b082fd36e8c1034af/e2e/$mdName"

$dateMain = "2016-09-06 10:59:25"
$dateZhHo = "2016-09-06 10:59:20"
$dateZhHb = "2016-09-06 10:59:39"
$dateDeHb = "2016-09-06 10:59:46"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 4
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

$wsO.Cells.Item(4, 1).Value = $mdName
$wsO.Cells.Item(4, 2).Value = $mdDisp
$wsO.Cells.Item(4, 3).Value = ".md"
$wsO.Cells.Item(4, 5).Value = $statusInSync
$wsO.Cells.Item(4, 6).Value = $statusInSync
$wsO.Cells.Item(4, 7).Value = $dateMain
$wsO.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsO.Hyperlinks.Add($wsO.Range("B4"), $mainUrl, "", "", $mdDisp)

$loO = $wsO.ListObjects.Item(1)
$loO.Resize($wsO.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4, 1).Value = $mdName
$wsZh.Cells.Item(4, 2).Value = ".md"
$wsZh.Cells.Item(4, 3).Value = $statusInSync
$wsZh.Cells.Item(4, 4).Value = "e2e"
$wsZh.Cells.Item(4, 5).Value = "ht"
$wsZh.Cells.Item(4, 6).Value = "True"
$wsZh.Cells.Item(4, 7).Value = $xlfZhCn
$wsZh.Cells.Item(4, 8).Value = $dateZhHo
$wsZh.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 9).Value = $mdName
$wsZh.Cells.Item(4, 10).Value = $xlfZhCn
$wsZh.Cells.Item(4, 11).Value = $dateZhHb
$wsZh.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 13).Value = "True"
$wsZh.Cells.Item(4, 15).Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mainUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhcnUrl, "", "", $mdName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4, 1).Value = $mdName
$wsDe.Cells.Item(4, 2).Value = ".md"
$wsDe.Cells.Item(4, 3).Value = $statusInSync
$wsDe.Cells.Item(4, 4).Value = "e2e"
$wsDe.Cells.Item(4, 5).Value = "ht"
$wsDe.Cells.Item(4, 6).Value = "True"
$wsDe.Cells.Item(4, 7).Value = $xlfDeDe
$wsDe.Cells.Item(4, 8).Value = $dateMain
$wsDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 9).Value = $mdName
$wsDe.Cells.Item(4, 10).Value = $xlfDeDe
$wsDe.Cells.Item(4, 11).Value = $dateDeHb
$wsDe.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 13).Value = "True"
$wsDe.Cells.Item(4, 15).Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mainUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $dedeUrl, "", "", $mdName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
